# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: bold header row for the new size-classification table
$ws.Range("B18").Value = "Number of employees"
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B18:D18").Style = "title"

# Rows 19-22 use the plain/default style (A21 & A22 previously held the
# bold "SNC" / italic source-citation cells, so their old formatting must
# be cleared back to the default now that they hold table data instead).
$ws.Range("A19:D22").ClearFormats()

# Row 19: Micro (B19/C19/D19 are explicit empty-text cells, like the source
# workbook - a leading apostrophe forces an empty text value rather than a
# truly blank cell)
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "'"
$ws.Range("C19").Value = "'"
$ws.Range("D19").Value = "'"

# Row 20: Small
$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "<50"
$ws.Range("C20").Value = "'"
$ws.Range("D20").Value = "<100,000 UT"

# Row 21: Medium
$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "51-100 <br/><250 Industry, <br/><500 Trade, <br/><100 Service, <br/><50 Agriculture"
$ws.Range("C21").Value = "'"
$ws.Range("D21").Value = "100,000 UT to 250,000 UT <br/><750,000 Industry, <br/><1,000,000 Trade, <br/><500,000 Serv., <br/><300,000 Agriculture"

# Row 22: Large
$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100 <br/>>=250 Industry, <br/>>=500 Trade, <br/>>=100 Service, <br/>>=50 Agriculture"
$ws.Range("C22").Value = "'"
$ws.Range("D22").Value = ">250,000 UT <br/>>=750,000 Industry, <br/>>=1,000,000 Trade, <br/>>=500,000 Serv.,<br/> >=300,000 Agriculture"

# Move the source citation (previously rows 21-22) down to rows 27-28
$ws.Range("A27").Value = "SNC"
$ws.Range("A27").Style = "title"
$ws.Range("A28").Value = "Servicio Nacional de Contrataciones (SNC), Foro de Estandares Internacionales, Papel del Contador Publico en la PYMEs Venezolanas. Available at http://fccpv.org/cont3/data/files/Foro-II-May2009-Presentacion-1.pdf"
$ws.Range("A28").Style = "source"
